# Logged Week 15 and simulated Week 16
# Update Target Depth Data for the Dolphins: "H" row (row 2) on both the
# OFF and DEF sheets gets updated totals after logging Week 15 and
# simulating Week 16.

$wb = $excel.ActiveWorkbook

# OFF sheet - row 2 (H)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 394
$wsOff.Range("C2").Value = 271
$wsOff.Range("D2").Value = 59
$wsOff.Range("E2").Value = 27
$wsOff.Range("F2").Value = 4
$wsOff.Range("G2").Value = 6

# DEF sheet - row 2 (H)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 504
$wsDef.Range("C2").Value = 336
$wsDef.Range("D2").Value = 136
$wsDef.Range("E2").Value = 60
